$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 7919
}

for ($r = 11; $r -le 22; $r++) {
    $ws.Cells.Item($r, 3).Value = 7917
}

for ($r = 139; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
